# This deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme"  (only ever referenced by the Notes Master)
#   ppt/theme/theme2.xml -> "Integral"      (referenced by the one-and-only Slide Master,
#                                             i.e. the theme that actually paints the slides)
#
# The authored change swaps the two themes' contents: the deck-visible theme
# (the Slide Master's) becomes the stock "Office Theme" colors, while the
# (previously) "Office Theme" colors move to the other slot. We reproduce the
# visible/deck-affecting half of that swap -- recoloring the active Slide
# Master theme from the "Integral" palette to the default "Office" palette --
# by rewriting every slot of the master's theme color scheme via the
# PowerPoint object model (SlideMaster.ColorScheme.Colors(i).RGB), which is
# the supported COM surface for editing a theme's 12 color slots
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).

$p  = $ppt.ActivePresentation
$sm = $p.SlideMaster
$cs = $sm.ColorScheme

function HexToBgr([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Office Theme color scheme, in the standard theme slot order:
# dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $cs.Colors($i).RGB = HexToBgr $officeThemeColors[$i - 1]
}
